# Add data for 2021-12-13 (commit message says 2021-12-21, but the diff
# content clearly updates dates from 12-12 to 12-13)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the sheet (tab) name to reflect the new "through" date.
$ws.Name = "Through 2021-12-13"

# Update the December row label text (column A, row 13).
$ws.Range("A13").Value = "December (through 12-13)"

# Update October 2021 value (row 11).
$ws.Range("H11").Value = 196

# Update December row (row 13) values for each year column.
$ws.Range("B13").Value = 13
$ws.Range("C13").Value = 38
$ws.Range("D13").Value = 44
$ws.Range("E13").Value = 28
$ws.Range("F13").Value = 22
$ws.Range("G13").Value = 66
$ws.Range("H13").Value = 95

# Update Total row (row 14) values for each year column.
$ws.Range("B14").Value = 304
$ws.Range("C14").Value = 601
$ws.Range("D14").Value = 865
$ws.Range("E14").Value = 710
$ws.Range("F14").Value = 556
$ws.Range("G14").Value = 1330
$ws.Range("H14").Value = 1738
